$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BECbIC")

$values = @{
    "B2"  = 131344000
    "C2"  = 575870000
    "D2"  = 369166000
    "E2"  = 1409375000
    "F2"  = 330422000
    "G2"  = 10088000
    "H2"  = 43596000
    "I2"  = 69451000
    "J2"  = 118685000
    "K2"  = 121941000
    "L2"  = 58553000
    "M2"  = 112433000
    "N2"  = 24283000
    "O2"  = 143330000
    "P2"  = 491016000
    "Q2"  = 40162000
    "R2"  = 96428000
    "S2"  = 0
    "T2"  = 0
    "U2"  = 99953000
    "V2"  = 500950000
    "W2"  = 2862120000
    "X2"  = 4727832000
    "Y2"  = 1445966000
    "Z2"  = 2125263000
    "AA2" = 370389000
    "AB2" = 383692000
    "AC2" = 42792000
    "AD2" = 1880567000
    "AE2" = 506937000
    "AF2" = 1293077000
    "AG2" = 14398860000
    "AH2" = 428925000
    "AI2" = 6262405000
    "AJ2" = 318637000
    "AK2" = 116370000
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
